# Weekly driver report update for 2025-04-19
# Refresh the "Good Drivers" table (rows 13-18) with this week's figures:
# client counts grow, a couple of drivers' "Good Roaming %" ticks up to
# 100, and the driver-vintage date is not yet known for newer drivers
# (those cells are left blank).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13: Intel(R) Wi-Fi 6 AX201 160MHz - 21.60.2.1
$ws.Range("A13").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.60.2.1"
$ws.Range("B13").Value = 56018
$ws.Range("D13").Value = 100
$ws.Range("E13").Value = ""

# Row 14: Intel(R) Wi-Fi 6 AX201 160MHz - 22.50.1.1
$ws.Range("A14").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.50.1.1"
$ws.Range("B14").Value = 34244
$ws.Range("D14").Value = 100
$ws.Range("E14").Value = ""

# Row 15: Intel(R) Wi-Fi 6 AX201 160MHz - 23.100.0.4
$ws.Range("A15").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 23.100.0.4"
$ws.Range("B15").Value = 442178
$ws.Range("D15").Value = 99.90000000000001
# Leading apostrophe forces text (not an auto-converted date serial), to
# match the vintage column's plain-text "YYYY-MM-DD" values.
$ws.Range("E15").Value = "'2024-11-10"

# Row 16: Intel(R) Wi-Fi 6 AX201 160MHz - 22.80.0.9
$ws.Range("A16").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.80.0.9"
$ws.Range("B16").Value = 77849
$ws.Range("D16").Value = 99.90000000000001
$ws.Range("E16").Value = "'2021-08-18"

# Row 17: Intel(R) Wi-Fi 6 AX201 160MHz - 21.110.3.2
$ws.Range("A17").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.110.3.2"
$ws.Range("B17").Value = 59673
$ws.Range("D17").Value = 100
$ws.Range("E17").Value = "'2020-08-05"

# Row 18: Intel(R) Wi-Fi 6 AX201 160MHz - 21.70.0.6
$ws.Range("A18").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.70.0.6"
$ws.Range("B18").Value = 113652
$ws.Range("D18").Value = 100
$ws.Range("E18").Value = ""
